{"js": "// Add a default header to the document's (first) section containing the\n// questionnaire number, so printed pages can be tracked.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst header = section.getHeader(\"Primary\");\n\nconst paragraph = header.insertParagraph(\"Questionnaire 15\", \"Replace\");\nparagraph.styleBuiltIn = Word.BuiltInStyleName.header;\nparagraph.alignment = Word.Alignment.centered;\nparagraph.font.name = \"Arial\";\nparagraph.font.size = 12;\n\nawait context.sync();\n", "ps1": "# Add a default (primary) header to the first/only section containing the\n# questionnaire number, so the questionnaire can be identified after printing.\n$d = $word.ActiveDocument\n$section = $d.Sections(1)\n$header = $section.Headers(1)  # wdHeaderFooterPrimary\n\n$header.Range.InsertAfter(\"Questionnaire 15\")\n$header.Range.Paragraphs(1).Style = \"Header\"\n$header.Range.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter\n\n$textRange = $header.Range.Paragraphs(1).Range\n[void]$textRange.MoveEnd(1, -1)  # exclude the paragraph mark from the font change\n$textRange.Font.Name = \"Arial\"\n$textRange.Font.Size = 12\n"}
